$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Ревизор - Жамшид", "2451", "ИП `"Худайбергенов`"", "Организация", "б/н", "76301", "1494", "1635", "141", "2025", "5", "23", "16", "19"),
    @("Ревизор - Жамшид", "41369", "Ип Омаров Ж", "Караб.шоссе", "б/н", "79214", "97", "116", "19", "2025", "5", "23", "16", "28"),
    @("Ревизор - Жамшид", "38290", "МОЙКА 24", "Организация", "б/н", "81248", "488", "673", "185", "2025", "5", "23", "16", "32"),
    @("Ревизор - Жамшид", "17634", "Ясли сад Кайнар", "Чкалова", "19.", "73687", "1546", "1568", "22", "2025", "5", "23", "16", "33"),
    @("Ревизор - Жамшид", "41623", "Ип Зияев", "Караб.шоссе", "б/н", "80496", "377", "417", "40", "2025", "5", "23", "16", "39")
)

$startRow = 70
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 1; $c -le $row.Length; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $row[$c - 1]
        if ($val -match '^-?[0-9]+(\.[0-9]*)?$') {
            # Values in this log are logically text (account numbers, dates
            # split into parts, etc.) even when they look numeric - keep
            # them stored as text like the rest of the sheet.
            $cell.NumberFormat = "@"
        }
        $cell.Value = $val
    }
}
